$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the row label text (shared string) for row 3
$ws.Range("A3").Value = "2020-06-29_diff"

# Update the numeric values in row 3
$ws.Range("B3").Value = -5.229973662
$ws.Range("C3").Value = 8.286826071

# D3 is no longer populated in the updated scope; clear its contents
$ws.Range("D3").ClearContents()
